$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '28.101.10'
Set-TextValue 'E2' '  +0.26%  '

Set-TextValue 'D3' '1.873.34'
Set-TextValue 'E3' '  +0.49%  '

Set-TextValue 'D4' '1.003'
Set-TextValue 'E4' '  -0.01%  '

Set-TextValue 'D5' '312.72'
Set-TextValue 'E5' '  +0.20%  '

Set-TextValue 'D6' '1.001'
Set-TextValue 'E6' '  -0.07%  '

Set-TextValue 'D7' '0.5109'
Set-TextValue 'E7' '  +0.32%  '

Set-TextValue 'D8' '0.3881'
Set-TextValue 'E8' '  +1.44%  '

Set-TextValue 'D9' '0.08352'
Set-TextValue 'E9' '  +0.49%  '

Set-TextValue 'D10' '1.116'
Set-TextValue 'E10' '  +0.08%  '

Set-TextValue 'B11' 'OKB'
Set-TextValue 'C11' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D11' '41.48'
Set-TextValue 'E11' '  -0.19%  '

Set-TextValue 'B12' 'Polkadot'
Set-TextValue 'C12' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D12' '6.228'
Set-TextValue 'E12' '  +0.11%  '

Set-TextValue 'B13' 'Solana'
Set-TextValue 'C13' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D13' '20.55'
Set-TextValue 'E13' '  -0.25%  '

Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.861.98'
Set-TextValue 'E14' '  +0.13%  '

Set-TextValue 'B15' 'Chainlink'
Set-TextValue 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '7.250'
Set-TextValue 'E15' '  +0.54%  '

Set-TextValue 'B16' 'BinanceUSD'
Set-TextValue 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D16' '1.002'
Set-TextValue 'E16' '  -0.05%  '

Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.00001101'
Set-TextValue 'E17' '  +0.36%  '

Set-TextValue 'B18' 'Litecoin'
Set-TextValue 'C18' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D18' '90.62'
Set-TextValue 'E18' '  -0.27%  '

Set-TextValue 'B19' 'TRON'
Set-TextValue 'C19' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D19' '0.06669'
Set-TextValue 'E19' '  +0.61%  '

Set-TextValue 'B20' 'Avalanche'
Set-TextValue 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '17.73'
Set-TextValue 'E20' '  +0.01%  '

Set-TextValue 'B21' 'Dai'
Set-TextValue 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D21' '1.001'
Set-TextValue 'E21' '  -0.13%  '

Set-TextValue 'B22' 'Uniswap'
Set-TextValue 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '6.003'
Set-TextValue 'E22' '  -0.55%  '

Set-TextValue 'B23' 'WrappedBTC'
Set-TextValue 'C23' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D23' '28.128.76'
Set-TextValue 'E23' '  +0.29%  '

Set-TextValue 'B24' 'Cosmos'
Set-TextValue 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D24' '11.10'
Set-TextValue 'E24' '  +0.23%  '

Set-TextValue 'B25' 'Toncoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.233'
Set-TextValue 'E25' '  -0.16%  '

Set-TextValue 'B26' 'Monero'
Set-TextValue 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '158.55'
Set-TextValue 'E26' '  +0.21%  '

Set-TextValue 'B27' 'LidoDAOToken'
Set-TextValue 'C27' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D27' '2.464'
Set-TextValue 'E27' '  -3.39%  '

Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '20.56'
Set-TextValue 'E28' '  -0.22%  '

Set-TextValue 'B29' 'BitcoinCash'
Set-TextValue 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D29' '124.86'
Set-TextValue 'E29' '  -0.39%  '

Set-TextValue 'B30' 'Stellar'
Set-TextValue 'C30' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D30' '0.1051'
Set-TextValue 'E30' '  -0.24%  '

Set-TextValue 'B31' 'ImmutableX'
Set-TextValue 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '1.034'
Set-TextValue 'E31' '  -0.51%  '

Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '5.804'
Set-TextValue 'E32' '  -0.86%  '

Set-TextValue 'B33' 'HuobiToken'
Set-TextValue 'C33' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D33' '3.591'
Set-TextValue 'E33' '  -0.18%  '

Set-TextValue 'B34' 'FraxShare'
Set-TextValue 'C34' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D34' '9.606'
Set-TextValue 'E34' '  +1.84%  '

Set-TextValue 'B35' 'VeChain'
Set-TextValue 'C35' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D35' '0.02440'
Set-TextValue 'E35' '  +1.06%  '

Set-TextValue 'B36' 'Hedera'
Set-TextValue 'C36' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D36' '0.06545'
Set-TextValue 'E36' '  +0.29%  '

Set-TextValue 'B37' 'Algorand'
Set-TextValue 'C37' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D37' '0.2200'
Set-TextValue 'E37' '  +1.43%  '

Set-TextValue 'B38' 'ARBITRUM'
Set-TextValue 'C38' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D38' '1.194'
Set-TextValue 'E38' '  -1.07%  '

Set-TextValue 'B39' 'TheSandbox'
Set-TextValue 'C39' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D39' '0.6473'
Set-TextValue 'E39' '  +0.09%  '

Set-TextValue 'B40' 'TrustWalletToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D40' '1.232'
Set-TextValue 'E40' '  +0.60%  '

Set-TextValue 'B41' 'InternetComputer(DFINITY)'
Set-TextValue 'C41' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D41' '4.956'
Set-TextValue 'E41' '  +0.37%  '

Set-TextValue 'B42' 'Aptos'
Set-TextValue 'C42' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D42' '11.24'
Set-TextValue 'E42' '  +0.14%  '

Set-TextValue 'B43' 'Decentraland'
Set-TextValue 'C43' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D43' '0.6094'
Set-TextValue 'E43' '  -0.21%  '

Set-TextValue 'B44' 'EnergySwap'
Set-TextValue 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D44' '12.96'
Set-TextValue 'E44' '  -1.21%  '

Set-TextValue 'B45' 'WEMIXTOKEN'
Set-TextValue 'C45' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D45' '1.276'
Set-TextValue 'E45' '  -0.63%  '

Set-TextValue 'B46' 'PancakeSwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D46' '3.657'
Set-TextValue 'E46' '  -0.39%  '

Set-TextValue 'B47' 'NEARProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D47' '2.005'
Set-TextValue 'E47' '  -0.62%  '

Set-TextValue 'B48' 'EOS'
Set-TextValue 'C48' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D48' '1.231'
Set-TextValue 'E48' '  +1.85%  '

Set-TextValue 'B49' 'Quant'
Set-TextValue 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D49' '120.39'
Set-TextValue 'E49' '  -0.08%  '

Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.06894'
Set-TextValue 'E50' '  +0.51%  '

Set-TextValue 'B51' 'Aave'
Set-TextValue 'C51' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D51' '77.60'
Set-TextValue 'E51' '  -1.17%  '
